$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# The document has a "first page" header/footer pair and a "default" (odd/even)
# header/footer pair, each containing one inline logo picture:
#   Headers.Item(1) -> header2.xml (default)  - BTec_Logo-Orange, id=3
#   Headers.Item(2) -> header1.xml (first)    - BTec_Logo-Orange, id=1
#   Footers.Item(1) -> footer2.xml (default)  - PearsonLogo,      id=4
#   Footers.Item(2) -> footer1.xml (first)    - PearsonLogo,      id=2
#
# Rename the embedded picture objects: the Pearson logos go from image1.png to
# image2.png, and the BTec logos go from image2.jpg to image1.jpg.

$hdr1 = $sec.Headers.Item(1).Range.InlineShapes.Item(1)
$hdr1.Name = "image1.jpg"

$hdr2 = $sec.Headers.Item(2).Range.InlineShapes.Item(1)
$hdr2.Name = "image1.jpg"

$ftr1 = $sec.Footers.Item(1).Range.InlineShapes.Item(1)
$ftr1.Name = "image2.png"

$ftr2 = $sec.Footers.Item(2).Range.InlineShapes.Item(1)
$ftr2.Name = "image2.png"
